$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - numeric values
$ws.Range("A2").Value = "14/07/2023"
$ws.Range("B2").Value = 6000
$ws.Range("C2").Value = 6000
$ws.Range("D2").Value = 6000
$ws.Range("E2").Value = 6000
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 100

# Row 3 - numeric values
$ws.Range("A3").Value = "14/07/2023"
$ws.Range("B3").Value = 5000
$ws.Range("C3").Value = 11000
$ws.Range("D3").Value = 6000
$ws.Range("E3").Value = 12000
$ws.Range("F3").Value = 20
$ws.Range("G3").Value = 1000
$ws.Range("H3").Value = 109.09

# Row 4 - numeric values
$ws.Range("A4").Value = "14/07/2023"
$ws.Range("B4").Value = 5000
$ws.Range("C4").Value = 16000
$ws.Range("D4").Value = 2000
$ws.Range("E4").Value = 14000
$ws.Range("F4").Value = 30
$ws.Range("G4").Value = 2000
$ws.Range("H4").Value = 87.5

# Row 5 - values entered as text (matching the new "R$ ..." backup formatting
# mentioned in the commit message) so they are stored as text, not numbers.
$ws.Range("A5").Value = "14/07/2023"
$ws.Range("B5").Value = "'5000.00"
$ws.Range("C5").Value = "'21000.00"
$ws.Range("D5").Value = "'7000.00"
$ws.Range("E5").Value = "'21000.00"
$ws.Range("F5").Value = "'40.0"
$ws.Range("G5").Value = "'0.00"
$ws.Range("H5").Value = "'100.00"
